# Update the "panel_query_time" timestamps on the "data" sheet (F2:F9)
# from the 10:49:56.xxxxxx run to the 14:33:03.xxxxxx run.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:33:03.266002"
$data.Range("F3").Value = "2021-10-05 14:33:03.266010"
$data.Range("F4").Value = "2021-10-05 14:33:03.266013"
$data.Range("F5").Value = "2021-10-05 14:33:03.266016"
$data.Range("F6").Value = "2021-10-05 14:33:03.266019"
$data.Range("F7").Value = "2021-10-05 14:33:03.266021"
$data.Range("F8").Value = "2021-10-05 14:33:03.266024"
$data.Range("F9").Value = "2021-10-05 14:33:03.266026"

# Add a new "metadata" worksheet right after "data" describing the panel query.
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the bold/bordered/centered header style from "data"!B1 onto the new headers.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Achromatopsia"
$meta.Range("C2").Value = 3149
$meta.Range("D2").Value = "'1.3"
$meta.Range("E2").Value = "2020-11-02T06:54:44.503816Z"
$meta.Range("F2").Value = "2021-10-05 14:33:03.261963"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3149/?format=json"

# Keep "data" as the active/selected sheet, matching the original workbook state.
$data.Activate()
